$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-RowValues($ws, $rowRef, $values) {
    $n = $values.Length
    $arr = New-Object 'object[,]' 1,$n
    for ($i = 0; $i -lt $n; $i++) {
        $arr[0,$i] = $values[$i]
    }
    $ws.Range($rowRef).Value = $arr
}

# Step 1: the diff drops the last column (old BA) and the last two rows (old 23, 24)
$ws.Range("BA:BA").Delete()
$ws.Range("23:24").Delete()

# Step 2: rewrite every surviving row (A1:AZ22) with the corrected values from the fix
Set-RowValues $ws "A1:AZ1" @($null, 39583, 39765, 39948, 40130, 40310, 40494, 40676, 40862, 41044, 41228, 41409, 41592, 41774, 41957, 42137, 42321, 42503, 42689, 42867, 43053, 43145, 43235, 43326, 43418, 43510, 43600, 43691, 43783, 43875, 43966, 44068, 44159, 44251, 44341, 44432, 44525, 44617, 44706, 44798, 44890, 44981, 45071, 45163, 45254, 45345, 45436, 45534, 45618, 45713, 45800, 45891)
Set-RowValues $ws "A2:AZ2" @(39813, $null, $null, $null, $null, $null, $null, $null, $null, $null, $null, $null, $null, $null, $null, $null, $null, $null, $null, $null, $null, $null, $null, $null, $null, $null, $null, $null, $null, $null, $null, $null, $null, $null, $null, $null, $null, $null, $null, $null, $null, $null, $null, $null, $null, $null, $null, $null, $null, $null, $null, $null)
Set-RowValues $ws "A3:AZ3" @(40178, 0.5087393606160395, -0.4513776153963867, -1.118515468742087, -1.324983933426882, -1.324983933426882, -1.324983933426882, -1.324983933426882, -1.324983933426882, -1.324983933426882, -1.324983933426882, -1.324983933426882, -1.324983933426882, -1.324983933426882, -1.324983933426882, -1.324983933426882, -1.324983933426882, -1.324983933426882, -1.324983933426882, -1.324983933426882, -1.324983933426882, -1.324983933426882, -1.324983933426882, -1.324983933426882, -1.324983933426882, -1.324983933426882, -1.324983933426882, -1.324983933426882, -1.324983933426882, -1.324983933426882, -1.324983933426882, -1.324983933426882, -1.324983933426882, -1.324983933426882, -1.324983933426882, -1.324983933426882, -1.324983933426882, -1.324983933426882, -1.324983933426882, -1.324983933426882, -1.324983933426882, -1.324983933426882, -1.324983933426882, -1.324983933426882, -1.324983933426882, -1.324983933426882, -1.324983933426882, -1.324983933426882, -1.324983933426882, -1.324983933426882, -1.324983933426882, -1.324983933426882)
Set-RowValues $ws "A4:AZ4" @(40543, 0.5613599890424092, -0.2663267917771495, -0.6296678961043134, -0.8803581938132576, -0.1156872058426073, -0.3900454704678369, -0.3900454704678369, -0.3900454704678369, -0.3900454704678369, -0.3900454704678369, -0.3900454704678369, -0.3900454704678369, -0.3900454704678369, -0.3900454704678369, -0.3900454704678369, -0.3900454704678369, -0.3900454704678369, -0.3900454704678369, -0.3900454704678369, -0.3900454704678369, -0.3900454704678369, -0.3900454704678369, -0.3900454704678369, -0.3900454704678369, -0.3900454704678369, -0.3900454704678369, -0.3900454704678369, -0.3900454704678369, -0.3900454704678369, -0.3900454704678369, -0.3900454704678369, -0.3900454704678369, -0.3900454704678369, -0.3900454704678369, -0.3900454704678369, -0.3900454704678369, -0.3900454704678369, -0.3900454704678369, -0.3900454704678369, -0.3900454704678369, -0.3900454704678369, -0.3900454704678369, -0.3900454704678369, -0.3900454704678369, -0.3900454704678369, -0.3900454704678369, -0.3900454704678369, -0.3900454704678369, -0.3900454704678369, -0.3900454704678369, -0.3900454704678369)
Set-RowValues $ws "A5:AZ5" @(40908, $null, $null, -0.6591481056594906, -0.8332792306556258, -0.5120992642018263, -1.213027585730386, -0.4084169314491404, -0.2995848153489522, -0.2995848153489522, -0.2995848153489522, -0.2995848153489522, -0.2995848153489522, -0.2995848153489522, -0.2995848153489522, -0.2995848153489522, -0.2995848153489522, -0.2995848153489522, -0.2995848153489522, -0.2995848153489522, -0.2995848153489522, -0.2995848153489522, -0.2995848153489522, -0.2995848153489522, -0.2995848153489522, -0.2995848153489522, -0.2995848153489522, -0.2995848153489522, -0.2995848153489522, -0.2995848153489522, -0.2995848153489522, -0.2995848153489522, -0.2995848153489522, -0.2995848153489522, -0.2995848153489522, -0.2995848153489522, -0.2995848153489522, -0.2995848153489522, -0.2995848153489522, -0.2995848153489522, -0.2995848153489522, -0.2995848153489522, -0.2995848153489522, -0.2995848153489522, -0.2995848153489522, -0.2995848153489522, -0.2995848153489522, -0.2995848153489522, -0.2995848153489522, -0.2995848153489522, -0.2995848153489522, -0.2995848153489522)
Set-RowValues $ws "A6:AZ6" @(41274, $null, $null, $null, $null, -0.5336548278322195, -1.019482592112642, -0.6403426624573716, -0.3230872999110068, -0.2188016966516937, -0.2075757021743008, -0.2075757021743008, -0.2075757021743008, -0.2075757021743008, -0.2075757021743008, -0.2075757021743008, -0.2075757021743008, -0.2075757021743008, -0.2075757021743008, -0.2075757021743008, -0.2075757021743008, -0.2075757021743008, -0.2075757021743008, -0.2075757021743008, -0.2075757021743008, -0.2075757021743008, -0.2075757021743008, -0.2075757021743008, -0.2075757021743008, -0.2075757021743008, -0.2075757021743008, -0.2075757021743008, -0.2075757021743008, -0.2075757021743008, -0.2075757021743008, -0.2075757021743008, -0.2075757021743008, -0.2075757021743008, -0.2075757021743008, -0.2075757021743008, -0.2075757021743008, -0.2075757021743008, -0.2075757021743008, -0.2075757021743008, -0.2075757021743008, -0.2075757021743008, -0.2075757021743008, -0.2075757021743008, -0.2075757021743008, -0.2075757021743008, -0.2075757021743008, -0.2075757021743008)
Set-RowValues $ws "A7:AZ7" @(41639, $null, $null, $null, $null, $null, $null, -0.6829735749296906, -0.4248987225036194, -0.1561757764150462, -0.2793004163246238, 0.05915234751026066, 0.124712275190686, 0.124712275190686, 0.124712275190686, 0.124712275190686, 0.124712275190686, 0.124712275190686, 0.124712275190686, 0.124712275190686, 0.124712275190686, 0.124712275190686, 0.124712275190686, 0.124712275190686, 0.124712275190686, 0.124712275190686, 0.124712275190686, 0.124712275190686, 0.124712275190686, 0.124712275190686, 0.124712275190686, 0.124712275190686, 0.124712275190686, 0.124712275190686, 0.124712275190686, 0.124712275190686, 0.124712275190686, 0.124712275190686, 0.124712275190686, 0.124712275190686, 0.124712275190686, 0.124712275190686, 0.124712275190686, 0.124712275190686, 0.124712275190686, 0.124712275190686, 0.124712275190686, 0.124712275190686, 0.124712275190686, 0.124712275190686, 0.124712275190686, 0.124712275190686)
Set-RowValues $ws "A8:AZ8" @(42004, $null, $null, $null, $null, $null, $null, $null, $null, -0.1931908880471678, -0.3668158601783045, 0.04624521867206965, -0.119752617912039, -0.2979029954603529, -0.255298189276465, -0.255298189276465, -0.255298189276465, -0.255298189276465, -0.255298189276465, -0.255298189276465, -0.255298189276465, -0.255298189276465, -0.255298189276465, -0.255298189276465, -0.255298189276465, -0.255298189276465, -0.255298189276465, -0.255298189276465, -0.255298189276465, -0.255298189276465, -0.255298189276465, -0.255298189276465, -0.255298189276465, -0.255298189276465, -0.255298189276465, -0.255298189276465, -0.255298189276465, -0.255298189276465, -0.255298189276465, -0.255298189276465, -0.255298189276465, -0.255298189276465, -0.255298189276465, -0.255298189276465, -0.255298189276465, -0.255298189276465, -0.255298189276465, -0.255298189276465, -0.255298189276465, -0.255298189276465, -0.255298189276465, -0.255298189276465)
Set-RowValues $ws "A9:AZ9" @(42369, $null, $null, $null, $null, $null, $null, $null, $null, $null, $null, 0.04932853730390541, -0.07066185312620821, -0.1124510725819206, -0.05946205208092747, 0.07317408757452348, 0.07418514192796266, 0.07418514192796266, 0.07418514192796266, 0.07418514192796266, 0.07418514192796266, 0.07418514192796266, 0.07418514192796266, 0.07418514192796266, 0.07418514192796266, 0.07418514192796266, 0.07418514192796266, 0.07418514192796266, 0.07418514192796266, 0.07418514192796266, 0.07418514192796266, 0.07418514192796266, 0.07418514192796266, 0.07418514192796266, 0.07418514192796266, 0.07418514192796266, 0.07418514192796266, 0.07418514192796266, 0.07418514192796266, 0.07418514192796266, 0.07418514192796266, 0.07418514192796266, 0.07418514192796266, 0.07418514192796266, 0.07418514192796266, 0.07418514192796266, 0.07418514192796266, 0.07418514192796266, 0.07418514192796266, 0.07418514192796266, 0.07418514192796266, 0.07418514192796266)
Set-RowValues $ws "A10:AZ10" @(42735, $null, $null, $null, $null, $null, $null, $null, $null, $null, $null, $null, $null, -0.1199460107991879, -0.07652469744219292, -0.002181547367274828, -0.001680662521774678, -0.06188089372189953, -0.07611406013281474, -0.07611406013281474, -0.07611406013281474, -0.07611406013281474, -0.07611406013281474, -0.07611406013281474, -0.07611406013281474, -0.07611406013281474, -0.07611406013281474, -0.07611406013281474, -0.07611406013281474, -0.07611406013281474, -0.07611406013281474, -0.07611406013281474, -0.07611406013281474, -0.07611406013281474, -0.07611406013281474, -0.07611406013281474, -0.07611406013281474, -0.07611406013281474, -0.07611406013281474, -0.07611406013281474, -0.07611406013281474, -0.07611406013281474, -0.07611406013281474, -0.07611406013281474, -0.07611406013281474, -0.07611406013281474, -0.07611406013281474, -0.07611406013281474, -0.07611406013281474, -0.07611406013281474, -0.07611406013281474, -0.07611406013281474)
Set-RowValues $ws "A11:AZ11" @(43100, $null, $null, $null, $null, $null, $null, $null, $null, $null, $null, $null, $null, $null, $null, -0.002407973313722778, -0.03430226911815826, -0.07932008107318644, -0.1247901924724348, -0.2199961235931358, -0.191300579729714, -0.191300579729714, -0.191300579729714, -0.191300579729714, -0.191300579729714, -0.191300579729714, -0.191300579729714, -0.191300579729714, -0.191300579729714, -0.191300579729714, -0.191300579729714, -0.191300579729714, -0.191300579729714, -0.191300579729714, -0.191300579729714, -0.191300579729714, -0.191300579729714, -0.191300579729714, -0.191300579729714, -0.191300579729714, -0.191300579729714, -0.191300579729714, -0.191300579729714, -0.191300579729714, -0.191300579729714, -0.191300579729714, -0.191300579729714, -0.191300579729714, -0.191300579729714, -0.191300579729714, -0.191300579729714, -0.191300579729714)
Set-RowValues $ws "A12:AZ12" @(43465, $null, $null, $null, $null, $null, $null, $null, $null, $null, $null, $null, $null, $null, $null, $null, $null, -0.06805249969644711, -0.03405723100528002, -0.1022879117640763, -0.05219951976568327, 0.0441865668729946, 0.0882025545300813, 0.1415113532986956, 0.0970330232288763, 0.0970330232288763, 0.0970330232288763, 0.0970330232288763, 0.0970330232288763, 0.0970330232288763, 0.0970330232288763, 0.0970330232288763, 0.0970330232288763, 0.0970330232288763, 0.0970330232288763, 0.0970330232288763, 0.0970330232288763, 0.0970330232288763, 0.0970330232288763, 0.0970330232288763, 0.0970330232288763, 0.0970330232288763, 0.0970330232288763, 0.0970330232288763, 0.0970330232288763, 0.0970330232288763, 0.0970330232288763, 0.0970330232288763, 0.0970330232288763, 0.0970330232288763, 0.0970330232288763, 0.0970330232288763)
Set-RowValues $ws "A13:AZ13" @(43830, $null, $null, $null, $null, $null, $null, $null, $null, $null, $null, $null, $null, $null, $null, $null, $null, $null, $null, -0.08624408439182885, -0.1021035481439614, -0.07120909843567613, -0.05573300569792217, 0.02926805735909976, -0.1345737582127748, -0.5756287392657988, -0.7844010209450802, -0.6919146680131605, -0.7407518902333265, -0.7407518902333265, -0.7407518902333265, -0.7407518902333265, -0.7407518902333265, -0.7407518902333265, -0.7407518902333265, -0.7407518902333265, -0.7407518902333265, -0.7407518902333265, -0.7407518902333265, -0.7407518902333265, -0.7407518902333265, -0.7407518902333265, -0.7407518902333265, -0.7407518902333265, -0.7407518902333265, -0.7407518902333265, -0.7407518902333265, -0.7407518902333265, -0.7407518902333265, -0.7407518902333265, -0.7407518902333265, -0.7407518902333265)
Set-RowValues $ws "A14:AZ14" @(44196, $null, $null, $null, $null, $null, $null, $null, $null, $null, $null, $null, $null, $null, $null, $null, $null, $null, $null, $null, $null, $null, -0.07026949721593567, -0.0279021829078685, -0.07686644652562524, -0.2394607875814136, -0.3786583343736716, -0.244860729922769, -0.4363737508290888, -0.7124953797697064, -0.1808804304865297, 0.3056679541520335, 0.3056679541520335, 0.3056679541520335, 0.3056679541520335, 0.3056679541520335, 0.3056679541520335, 0.3056679541520335, 0.3056679541520335, 0.3056679541520335, 0.3056679541520335, 0.3056679541520335, 0.3056679541520335, 0.3056679541520335, 0.3056679541520335, 0.3056679541520335, 0.3056679541520335, 0.3056679541520335, 0.3056679541520335, 0.3056679541520335, 0.3056679541520335, 0.3056679541520335)
Set-RowValues $ws "A15:AZ15" @(44561, $null, $null, $null, $null, $null, $null, $null, $null, $null, $null, $null, $null, $null, $null, $null, $null, $null, $null, $null, $null, $null, $null, $null, $null, $null, -0.3185394740402248, -0.2641947304525827, -0.334882695390748, -0.4617192974095352, -0.1077309791980285, 0.3245880452514394, -0.514812792200714, -0.8680533514735522, -0.8769761459347714, -1.388491535160907, -1.388491535160907, -1.388491535160907, -1.388491535160907, -1.388491535160907, -1.388491535160907, -1.388491535160907, -1.388491535160907, -1.388491535160907, -1.388491535160907, -1.388491535160907, -1.388491535160907, -1.388491535160907, -1.388491535160907, -1.388491535160907, -1.388491535160907, -1.388491535160907)
Set-RowValues $ws "A16:AZ16" @(44926, $null, $null, $null, $null, $null, $null, $null, $null, $null, $null, $null, $null, $null, $null, $null, $null, $null, $null, $null, $null, $null, $null, $null, $null, $null, $null, $null, $null, $null, -0.2598822251329125, -0.03770935223248451, -0.3609163387359504, -0.5995895195426981, -0.5354669478056073, -2.541003699199929, -2.321721165370549, -1.867377038014506, -1.719168896439693, -1.678482969789596, -1.678482969789596, -1.678482969789596, -1.678482969789596, -1.678482969789596, -1.678482969789596, -1.678482969789596, -1.678482969789596, -1.678482969789596, -1.678482969789596, -1.678482969789596, -1.678482969789596, -1.678482969789596)
Set-RowValues $ws "A17:AZ17" @(45291, $null, $null, $null, $null, $null, $null, $null, $null, $null, $null, $null, $null, $null, $null, $null, $null, $null, $null, $null, $null, $null, $null, $null, $null, $null, $null, $null, $null, $null, $null, $null, $null, -0.498361054968377, -0.4007230315998478, -1.135765464024774, -1.191207469356215, -0.8537083312609495, -0.5989817782328322, -0.6229862770763095, -1.107351089172237, -0.9016470784766528, -0.6902657121583777, -0.5999457276250508, -0.5999457276250508, -0.5999457276250508, -0.5999457276250508, -0.5999457276250508, -0.5999457276250508, -0.5999457276250508, -0.5999457276250508, -0.5999457276250508)
Set-RowValues $ws "A18:AZ18" @(45657, $null, $null, $null, $null, $null, $null, $null, $null, $null, $null, $null, $null, $null, $null, $null, $null, $null, $null, $null, $null, $null, $null, $null, $null, $null, $null, $null, $null, $null, $null, $null, $null, $null, $null, $null, $null, -0.9810195867723248, -0.8414662684835927, -1.076812708919683, -1.716216309145657, -1.556833564431637, -1.134712300966823, -0.7253995615808195, -0.4628630633218611, -0.101315145211045, -0.07874066250703748, -0.05499271238530445, -0.05499271238530445, -0.05499271238530445, -0.05499271238530445, -0.05499271238530445)
Set-RowValues $ws "A19:AZ19" @(46022, $null, $null, $null, $null, $null, $null, $null, $null, $null, $null, $null, $null, $null, $null, $null, $null, $null, $null, $null, $null, $null, $null, $null, $null, $null, $null, $null, $null, $null, $null, $null, $null, $null, $null, $null, $null, $null, $null, $null, $null, -1.612451964998651, -1.331281919762917, -1.109950047731278, -1.061602237248327, -0.1859728711464226, -0.1533081888441812, -0.1365897193907339, -0.01934819856548309, 0.4236015715998187, 0.4305325812036687, 0.4335297397760618)
Set-RowValues $ws "A20:AZ20" @(46387, $null, $null, $null, $null, $null, $null, $null, $null, $null, $null, $null, $null, $null, $null, $null, $null, $null, $null, $null, $null, $null, $null, $null, $null, $null, $null, $null, $null, $null, $null, $null, $null, $null, $null, $null, $null, $null, $null, $null, $null, $null, $null, $null, $null, -0.1983696069224883, -0.1849084361995157, -0.2720151998323761, -0.2680812993602677, 0.05974827491124213, 0.2136583044595852, 0.2794570629465865)
Set-RowValues $ws "A21:AZ21" @(46752, $null, $null, $null, $null, $null, $null, $null, $null, $null, $null, $null, $null, $null, $null, $null, $null, $null, $null, $null, $null, $null, $null, $null, $null, $null, $null, $null, $null, $null, $null, $null, $null, $null, $null, $null, $null, $null, $null, $null, $null, $null, $null, $null, $null, $null, $null, $null, $null, -0.0614876679608467, 0.1343289496879674, 0.2251918092589511)
Set-RowValues $ws "A22:AZ22" @(47118, $null, $null, $null, $null, $null, $null, $null, $null, $null, $null, $null, $null, $null, $null, $null, $null, $null, $null, $null, $null, $null, $null, $null, $null, $null, $null, $null, $null, $null, $null, $null, $null, $null, $null, $null, $null, $null, $null, $null, $null, $null, $null, $null, $null, $null, $null, $null, $null, $null, $null, $null)
